$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original style of column D (rows 2-51) so that forcing text format
# for numeric-looking values does not leave a lasting style change on the cells.
$dRange = $ws.Range("D2:D51")
$dOrigStyle = $dRange.Style
$dRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "29.459.71"
$ws.Cells.Item(2, 5).Value = "  -3.03%  "
$ws.Cells.Item(3, 4).Value = "1.989.26"
$ws.Cells.Item(3, 5).Value = "  -5.00%  "
$ws.Cells.Item(4, 4).Value = "1.018"
$ws.Cells.Item(4, 5).Value = "  +1.54%  "
$ws.Cells.Item(5, 4).Value = "328.04"
$ws.Cells.Item(5, 5).Value = "  -4.29%  "
$ws.Cells.Item(6, 4).Value = "1.015"
$ws.Cells.Item(6, 5).Value = "  +1.34%  "
$ws.Cells.Item(7, 4).Value = "0.4971"
$ws.Cells.Item(7, 5).Value = "  -4.93%  "
$ws.Cells.Item(8, 4).Value = "0.4195"
$ws.Cells.Item(8, 5).Value = "  -5.11%  "
$ws.Cells.Item(9, 4).Value = "53.74"
$ws.Cells.Item(9, 5).Value = "  -1.55%  "
$ws.Cells.Item(10, 4).Value = "0.08885"
$ws.Cells.Item(10, 5).Value = "  -4.74%  "
$ws.Cells.Item(11, 4).Value = "1.106"
$ws.Cells.Item(11, 5).Value = "  -5.30%  "
$ws.Cells.Item(12, 4).Value = "23.10"
$ws.Cells.Item(12, 5).Value = "  -7.06%  "
$ws.Cells.Item(13, 4).Value = "1.993.05"
$ws.Cells.Item(13, 5).Value = "  -1.88%  "
$ws.Cells.Item(14, 4).Value = "7.912"
$ws.Cells.Item(14, 5).Value = "  -7.78%  "
$ws.Cells.Item(15, 4).Value = "6.405"
$ws.Cells.Item(15, 5).Value = "  -7.18%  "
$ws.Cells.Item(16, 4).Value = "1.017"
$ws.Cells.Item(16, 5).Value = "  +1.48%  "
$ws.Cells.Item(17, 4).Value = "93.32"
$ws.Cells.Item(17, 5).Value = "  -7.94%  "
$ws.Cells.Item(18, 4).Value = "0.00001101"
$ws.Cells.Item(18, 5).Value = "  -4.98%  "
$ws.Cells.Item(19, 4).Value = "0.06657"
$ws.Cells.Item(19, 5).Value = "  -0.11%  "
$ws.Cells.Item(20, 4).Value = "19.34"
$ws.Cells.Item(20, 5).Value = "  -8.49%  "
$ws.Cells.Item(21, 4).Value = "1.013"
$ws.Cells.Item(21, 5).Value = "  +1.24%  "
$ws.Cells.Item(22, 4).Value = "5.883"
$ws.Cells.Item(22, 5).Value = "  -7.04%  "
$ws.Cells.Item(23, 4).Value = "29.525.68"
$ws.Cells.Item(23, 5).Value = "  -2.84%  "
$ws.Cells.Item(24, 4).Value = "11.90"
$ws.Cells.Item(24, 5).Value = "  -5.09%  "
$ws.Cells.Item(25, 4).Value = "2.300"
$ws.Cells.Item(25, 5).Value = "  -0.23%  "
$ws.Cells.Item(26, 4).Value = "157.38"
$ws.Cells.Item(26, 5).Value = "  -3.47%  "
$ws.Cells.Item(27, 4).Value = "20.59"
$ws.Cells.Item(27, 5).Value = "  -5.68%  "
$ws.Cells.Item(28, 4).Value = "6.232"
$ws.Cells.Item(28, 5).Value = "  -8.75%  "
$ws.Cells.Item(29, 4).Value = "2.277"
$ws.Cells.Item(29, 5).Value = "  -9.14%  "
$ws.Cells.Item(30, 4).Value = "126.89"
$ws.Cells.Item(30, 5).Value = "  -4.72%  "
$ws.Cells.Item(31, 4).Value = "1.042"
$ws.Cells.Item(31, 5).Value = "  -8.36%  "
$ws.Cells.Item(32, 4).Value = "0.09863"
$ws.Cells.Item(32, 5).Value = "  -5.74%  "
$ws.Cells.Item(33, 4).Value = "1.541"
$ws.Cells.Item(33, 5).Value = "  -7.13%  "
$ws.Cells.Item(34, 4).Value = "3.806"
$ws.Cells.Item(34, 5).Value = "  -1.19%  "
$ws.Cells.Item(35, 4).Value = "5.774"
$ws.Cells.Item(35, 5).Value = "  -7.69%  "
$ws.Cells.Item(36, 4).Value = "0.02442"
$ws.Cells.Item(36, 5).Value = "  -7.33%  "
$ws.Cells.Item(37, 4).Value = "9.204"
$ws.Cells.Item(37, 5).Value = "  -9.15%  "
$ws.Cells.Item(38, 4).Value = "1.294"
$ws.Cells.Item(38, 5).Value = "  -3.32%  "
$ws.Cells.Item(39, 4).Value = "0.06335"
$ws.Cells.Item(39, 5).Value = "  -7.31%  "
$ws.Cells.Item(40, 4).Value = "0.6482"
$ws.Cells.Item(40, 5).Value = "  -7.14%  "
$ws.Cells.Item(41, 4).Value = "11.50"
$ws.Cells.Item(41, 5).Value = "  -8.54%  "
$ws.Cells.Item(42, 4).Value = "0.2028"
$ws.Cells.Item(42, 5).Value = "  -8.37%  "
$ws.Cells.Item(43, 5).Value = "  +1.19%  "
$ws.Cells.Item(44, 4).Value = "0.6265"
$ws.Cells.Item(44, 5).Value = "  -7.96%  "
$ws.Cells.Item(45, 4).Value = "13.43"
$ws.Cells.Item(45, 5).Value = "  -6.36%  "
$ws.Cells.Item(46, 4).Value = "2.180"
$ws.Cells.Item(46, 5).Value = "  -7.08%  "
$ws.Cells.Item(47, 4).Value = "1.302"
$ws.Cells.Item(47, 5).Value = "  -5.22%  "
$ws.Cells.Item(48, 4).Value = "3.497"
$ws.Cells.Item(48, 5).Value = "  -3.73%  "
$ws.Cells.Item(49, 4).Value = "0.00000000338"
$ws.Cells.Item(49, 5).Value = "  -1.14%  "
$ws.Cells.Item(50, 4).Value = "0.06950"
$ws.Cells.Item(50, 5).Value = "  -4.02%  "
$ws.Cells.Item(51, 4).Value = "1.117"
$ws.Cells.Item(51, 5).Value = "  -9.34%  "

# Restore original style/number format for column D
$dRange.Style = $dOrigStyle
